{"js": "// Update the date title and every two-digit multiplication answer cell\n// in document order. The document is: 1 title paragraph (\"YYYY-MM-DD Weekday\")\n// followed by a 5-column table whose cells each hold exactly one paragraph\n// of the form \"A\u00d7B=C\". The mapping below pairs each paragraph's current\n// (old) text with its replacement (new) text, in body order.\nconst REPLACEMENTS = [[\"2023-03-07 Tuesday\", \"2023-03-08 Wednesday\"], [\"23\u00d773=1679\", \"14\u00d746=644\"], [\"69\u00d785=5865\", \"66\u00d786=5676\"], [\"32\u00d778=2496\", \"81\u00d711=891\"], [\"18\u00d742=756\", \"20\u00d775=1500\"], [\"11\u00d730=330\", \"61\u00d780=4880\"], [\"99\u00d720=1980\", \"96\u00d770=6720\"], [\"42\u00d725=1050\", \"20\u00d733=660\"], [\"59\u00d726=1534\", \"48\u00d764=3072\"], [\"82\u00d778=6396\", \"39\u00d794=3666\"], [\"19\u00d797=1843\", \"87\u00d733=2871\"], [\"18\u00d775=1350\", \"84\u00d760=5040\"], [\"44\u00d739=1716\", \"57\u00d748=2736\"], [\"57\u00d718=1026\", \"70\u00d759=4130\"], [\"30\u00d772=2160\", \"80\u00d766=5280\"], [\"52\u00d721=1092\", \"19\u00d763=1197\"], [\"19\u00d767=1273\", \"21\u00d761=1281\"], [\"95\u00d717=1615\", \"48\u00d739=1872\"], [\"60\u00d799=5940\", \"78\u00d773=5694\"], [\"99\u00d793=9207\", \"10\u00d715=150\"], [\"53\u00d746=2438\", \"56\u00d728=1568\"], [\"24\u00d722=528\", \"57\u00d724=1368\"], [\"82\u00d712=984\", \"45\u00d794=4230\"], [\"14\u00d731=434\", \"84\u00d721=1764\"], [\"39\u00d751=1989\", \"35\u00d725=875\"], [\"96\u00d724=2304\", \"20\u00d715=300\"], [\"82\u00d775=6150\", \"91\u00d752=4732\"], [\"27\u00d728=756\", \"34\u00d732=1088\"], [\"83\u00d794=7802\", \"98\u00d713=1274\"], [\"62\u00d749=3038\", \"90\u00d733=2970\"], [\"61\u00d734=2074\", \"87\u00d738=3306\"], [\"70\u00d716=1120\", \"47\u00d760=2820\"], [\"97\u00d728=2716\", \"76\u00d796=7296\"], [\"84\u00d745=3780\", \"62\u00d715=930\"], [\"75\u00d776=5700\", \"77\u00d763=4851\"], [\"72\u00d741=2952\", \"96\u00d786=8256\"], [\"50\u00d792=4600\", \"45\u00d734=1530\"], [\"77\u00d774=5698\", \"69\u00d753=3657\"], [\"20\u00d740=800\", \"48\u00d752=2496\"], [\"73\u00d712=876\", \"95\u00d785=8075\"], [\"71\u00d789=6319\", \"27\u00d735=945\"], [\"71\u00d793=6603\", \"26\u00d747=1222\"], [\"44\u00d734=1496\", \"81\u00d789=7209\"], [\"22\u00d714=308\", \"95\u00d714=1330\"], [\"54\u00d717=918\", \"86\u00d779=6794\"], [\"75\u00d728=2100\", \"47\u00d785=3995\"], [\"72\u00d764=4608\", \"91\u00d724=2184\"], [\"25\u00d755=1375\", \"62\u00d724=1488\"], [\"72\u00d767=4824\", \"58\u00d740=2320\"], [\"12\u00d749=588\", \"55\u00d710=550\"], [\"60\u00d713=780\", \"45\u00d790=4050\"], [\"26\u00d777=2002\", \"32\u00d735=1120\"], [\"100\u00d782=8200\", \"57\u00d716=912\"], [\"12\u00d720=240\", \"89\u00d719=1691\"], [\"81\u00d744=3564\", \"41\u00d752=2132\"], [\"87\u00d774=6438\", \"12\u00d713=156\"], [\"39\u00d792=3588\", \"61\u00d750=3050\"], [\"62\u00d774=4588\", \"95\u00d733=3135\"], [\"21\u00d7100=2100\", \"53\u00d775=3975\"], [\"38\u00d744=1672\", \"47\u00d723=1081\"], [\"47\u00d774=3478\", \"10\u00d736=360\"], [\"69\u00d755=3795\", \"15\u00d772=1080\"], [\"50\u00d740=2000\", \"37\u00d726=962\"], [\"43\u00d785=3655\", \"87\u00d743=3741\"], [\"92\u00d724=2208\", \"69\u00d720=1380\"], [\"37\u00d784=3108\", \"79\u00d784=6636\"], [\"55\u00d750=2750\", \"28\u00d763=1764\"], [\"62\u00d713=806\", \"72\u00d780=5760\"], [\"45\u00d716=720\", \"86\u00d784=7224\"], [\"57\u00d713=741\", \"60\u00d767=4020\"], [\"36\u00d786=3096\", \"17\u00d721=357\"], [\"65\u00d748=3120\", \"46\u00d793=4278\"], [\"70\u00d751=3570\", \"36\u00d765=2340\"], [\"65\u00d743=2795\", \"100\u00d786=8600\"], [\"20\u00d752=1040\", \"52\u00d752=2704\"], [\"73\u00d780=5840\", \"68\u00d740=2720\"], [\"19\u00d759=1121\", \"80\u00d738=3040\"], [\"45\u00d773=3285\", \"30\u00d728=840\"], [\"22\u00d738=836\", \"63\u00d785=5355\"], [\"14\u00d771=994\", \"87\u00d731=2697\"], [\"63\u00d746=2898\", \"31\u00d710=310\"], [\"99\u00d782=8118\", \"85\u00d752=4420\"], [\"97\u00d734=3298\", \"88\u00d737=3256\"], [\"12\u00d736=432\", \"28\u00d720=560\"], [\"53\u00d795=5035\", \"32\u00d748=1536\"], [\"70\u00d794=6580\", \"56\u00d754=3024\"], [\"22\u00d719=418\", \"29\u00d725=725\"], [\"55\u00d795=5225\", \"17\u00d770=1190\"], [\"97\u00d745=4365\", \"25\u00d740=1000\"], [\"88\u00d789=7832\", \"65\u00d779=5135\"], [\"32\u00d771=2272\", \"50\u00d733=1650\"], [\"31\u00d790=2790\", \"67\u00d748=3216\"], [\"39\u00d759=2301\", \"39\u00d728=1092\"], [\"96\u00d761=5856\", \"62\u00d771=4402\"], [\"69\u00d748=3312\", \"37\u00d785=3145\"], [\"83\u00d781=6723\", \"71\u00d713=923\"], [\"11\u00d761=671\", \"80\u00d797=7760\"], [\"70\u00d731=2170\", \"95\u00d713=1235\"], [\"23\u00d777=1771\", \"33\u00d767=2211\"], [\"21\u00d759=1239\", \"80\u00d775=6000\"], [\"33\u00d785=2805\", \"24\u00d734=816\"]];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    `Expected ${REPLACEMENTS.length} paragraphs, found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = items[i];\n  // Sanity-check we are editing the paragraph the diff expects before\n  // overwriting it (guards against any unexpected structural drift).\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i} text mismatch: expected \"${oldText}\", found \"${para.text}\"`\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and every two-digit multiplication answer cell.\n# The document is a title paragraph followed by a 20x5 table; each cell\n# holds exactly one paragraph of the form \"A\u00d7B=C\". Values below mirror the\n# body order of the source diff (row-major, left to right, top to bottom).\n\n$d = $word.ActiveDocument\n\n$titleOld = \"2023-03-07 Tuesday\"\n$titleNew = \"2023-03-08 Wednesday\"\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text -ne ($titleOld + [char]13)) {\n    throw \"Title paragraph text mismatch: expected $titleOld, found $($titlePara.Range.Text)\"\n}\n$titlePara.Range.Text = $titleNew\n\n$answers = @(\n    @(\"23\u00d773=1679\", \"14\u00d746=644\"), @(\"69\u00d785=5865\", \"66\u00d786=5676\"), @(\"32\u00d778=2496\", \"81\u00d711=891\"), @(\"18\u00d742=756\", \"20\u00d775=1500\"), @(\"11\u00d730=330\", \"61\u00d780=4880\"),\n    @(\"99\u00d720=1980\", \"96\u00d770=6720\"), @(\"42\u00d725=1050\", \"20\u00d733=660\"), @(\"59\u00d726=1534\", \"48\u00d764=3072\"), @(\"82\u00d778=6396\", \"39\u00d794=3666\"), @(\"19\u00d797=1843\", \"87\u00d733=2871\"),\n    @(\"18\u00d775=1350\", \"84\u00d760=5040\"), @(\"44\u00d739=1716\", \"57\u00d748=2736\"), @(\"57\u00d718=1026\", \"70\u00d759=4130\"), @(\"30\u00d772=2160\", \"80\u00d766=5280\"), @(\"52\u00d721=1092\", \"19\u00d763=1197\"),\n    @(\"19\u00d767=1273\", \"21\u00d761=1281\"), @(\"95\u00d717=1615\", \"48\u00d739=1872\"), @(\"60\u00d799=5940\", \"78\u00d773=5694\"), @(\"99\u00d793=9207\", \"10\u00d715=150\"), @(\"53\u00d746=2438\", \"56\u00d728=1568\"),\n    @(\"24\u00d722=528\", \"57\u00d724=1368\"), @(\"82\u00d712=984\", \"45\u00d794=4230\"), @(\"14\u00d731=434\", \"84\u00d721=1764\"), @(\"39\u00d751=1989\", \"35\u00d725=875\"), @(\"96\u00d724=2304\", \"20\u00d715=300\"),\n    @(\"82\u00d775=6150\", \"91\u00d752=4732\"), @(\"27\u00d728=756\", \"34\u00d732=1088\"), @(\"83\u00d794=7802\", \"98\u00d713=1274\"), @(\"62\u00d749=3038\", \"90\u00d733=2970\"), @(\"61\u00d734=2074\", \"87\u00d738=3306\"),\n    @(\"70\u00d716=1120\", \"47\u00d760=2820\"), @(\"97\u00d728=2716\", \"76\u00d796=7296\"), @(\"84\u00d745=3780\", \"62\u00d715=930\"), @(\"75\u00d776=5700\", \"77\u00d763=4851\"), @(\"72\u00d741=2952\", \"96\u00d786=8256\"),\n    @(\"50\u00d792=4600\", \"45\u00d734=1530\"), @(\"77\u00d774=5698\", \"69\u00d753=3657\"), @(\"20\u00d740=800\", \"48\u00d752=2496\"), @(\"73\u00d712=876\", \"95\u00d785=8075\"), @(\"71\u00d789=6319\", \"27\u00d735=945\"),\n    @(\"71\u00d793=6603\", \"26\u00d747=1222\"), @(\"44\u00d734=1496\", \"81\u00d789=7209\"), @(\"22\u00d714=308\", \"95\u00d714=1330\"), @(\"54\u00d717=918\", \"86\u00d779=6794\"), @(\"75\u00d728=2100\", \"47\u00d785=3995\"),\n    @(\"72\u00d764=4608\", \"91\u00d724=2184\"), @(\"25\u00d755=1375\", \"62\u00d724=1488\"), @(\"72\u00d767=4824\", \"58\u00d740=2320\"), @(\"12\u00d749=588\", \"55\u00d710=550\"), @(\"60\u00d713=780\", \"45\u00d790=4050\"),\n    @(\"26\u00d777=2002\", \"32\u00d735=1120\"), @(\"100\u00d782=8200\", \"57\u00d716=912\"), @(\"12\u00d720=240\", \"89\u00d719=1691\"), @(\"81\u00d744=3564\", \"41\u00d752=2132\"), @(\"87\u00d774=6438\", \"12\u00d713=156\"),\n    @(\"39\u00d792=3588\", \"61\u00d750=3050\"), @(\"62\u00d774=4588\", \"95\u00d733=3135\"), @(\"21\u00d7100=2100\", \"53\u00d775=3975\"), @(\"38\u00d744=1672\", \"47\u00d723=1081\"), @(\"47\u00d774=3478\", \"10\u00d736=360\"),\n    @(\"69\u00d755=3795\", \"15\u00d772=1080\"), @(\"50\u00d740=2000\", \"37\u00d726=962\"), @(\"43\u00d785=3655\", \"87\u00d743=3741\"), @(\"92\u00d724=2208\", \"69\u00d720=1380\"), @(\"37\u00d784=3108\", \"79\u00d784=6636\"),\n    @(\"55\u00d750=2750\", \"28\u00d763=1764\"), @(\"62\u00d713=806\", \"72\u00d780=5760\"), @(\"45\u00d716=720\", \"86\u00d784=7224\"), @(\"57\u00d713=741\", \"60\u00d767=4020\"), @(\"36\u00d786=3096\", \"17\u00d721=357\"),\n    @(\"65\u00d748=3120\", \"46\u00d793=4278\"), @(\"70\u00d751=3570\", \"36\u00d765=2340\"), @(\"65\u00d743=2795\", \"100\u00d786=8600\"), @(\"20\u00d752=1040\", \"52\u00d752=2704\"), @(\"73\u00d780=5840\", \"68\u00d740=2720\"),\n    @(\"19\u00d759=1121\", \"80\u00d738=3040\"), @(\"45\u00d773=3285\", \"30\u00d728=840\"), @(\"22\u00d738=836\", \"63\u00d785=5355\"), @(\"14\u00d771=994\", \"87\u00d731=2697\"), @(\"63\u00d746=2898\", \"31\u00d710=310\"),\n    @(\"99\u00d782=8118\", \"85\u00d752=4420\"), @(\"97\u00d734=3298\", \"88\u00d737=3256\"), @(\"12\u00d736=432\", \"28\u00d720=560\"), @(\"53\u00d795=5035\", \"32\u00d748=1536\"), @(\"70\u00d794=6580\", \"56\u00d754=3024\"),\n    @(\"22\u00d719=418\", \"29\u00d725=725\"), @(\"55\u00d795=5225\", \"17\u00d770=1190\"), @(\"97\u00d745=4365\", \"25\u00d740=1000\"), @(\"88\u00d789=7832\", \"65\u00d779=5135\"), @(\"32\u00d771=2272\", \"50\u00d733=1650\"),\n    @(\"31\u00d790=2790\", \"67\u00d748=3216\"), @(\"39\u00d759=2301\", \"39\u00d728=1092\"), @(\"96\u00d761=5856\", \"62\u00d771=4402\"), @(\"69\u00d748=3312\", \"37\u00d785=3145\"), @(\"83\u00d781=6723\", \"71\u00d713=923\"),\n    @(\"11\u00d761=671\", \"80\u00d797=7760\"), @(\"70\u00d731=2170\", \"95\u00d713=1235\"), @(\"23\u00d777=1771\", \"33\u00d767=2211\"), @(\"21\u00d759=1239\", \"80\u00d775=6000\"), @(\"33\u00d785=2805\", \"24\u00d734=816\")\n)\n\n$tbl = $d.Tables.Item(1)\n$cols = $tbl.Columns.Count\n$rows = $tbl.Rows.Count\nif ($rows -ne 20 -or $cols -ne 5) {\n    throw \"Unexpected table shape: $rows x $cols\"\n}\n\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $idx = ($r - 1) * $cols + ($c - 1)\n        $pair = $answers[$idx]\n        $old = $pair[0]\n        $new = $pair[1]\n        $cell = $tbl.Cell($r, $c)\n        if ($cell.Range.Text -ne ($old + [char]13 + [char]7)) {\n            throw \"Cell ($r,$c) text mismatch: expected $old, found $($cell.Range.Text)\"\n        }\n        $cell.Range.Text = $new\n    }\n}\n\n"}
